# Weekly price update: insert a new "Jengibre" price record as the latest
# entry (row 40), pushing the existing historical rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 40, shifting rows 40..118 down to 41..119.
$ws.Rows.Item(40).Insert()

# Populate the new row with the latest observation.
$ws.Range("A40").Value = 6
$ws.Range("B40").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C40").Value = "Metropolitana"
$ws.Range("D40").Value = 45012
$ws.Range("E40").Value = 13
$ws.Range("F40").Value = 100114007
$ws.Range("G40").Value = "Jengibre"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 580
$ws.Range("K40").Value = 14000
$ws.Range("L40").Value = 15000
$ws.Range("M40").Value = 14448
$ws.Range("N40").Value = "$/caja 13 kilos"
$ws.Range("O40").Value = "Perú"
$ws.Range("P40").Value = 1111
$ws.Range("Q40").Value = 13
$ws.Range("R40").Value = "Hortaliza"
